$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the column whose header (row 15) contains "culture_collection" and delete that entire column.
$headerRow = 15
$lastCol = $ws.Cells.Item($headerRow, 201)
$col = $ws.Rows.Item($headerRow).Find("culture_collection")
if ($col -ne $null) {
    $col.EntireColumn.Delete()
}
